{"js": "// Revert \"para five added\": remove the last paragraph\n// (\"Para 5 \u2013 I love playing football and cricket.\") that was appended\n// to the document body.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst lastIndex = paragraphs.items.length - 1;\nconst lastParagraph = paragraphs.items[lastIndex];\n\nif (lastParagraph.text.indexOf(\"Para 5\") !== -1) {\n  lastParagraph.delete();\n  await context.sync();\n}\n", "ps1": "# Revert \"para five added\": remove the last paragraph\n# (\"Para 5 - I love playing football and cricket.\") that was appended\n# to the document body.\n\n$d = $word.ActiveDocument\n\n$n = $d.Paragraphs.Count\n$last = $d.Paragraphs($n)\n\nif ($last.Range.Text -match \"Para 5\") {\n    $last.Range.Delete()\n}\n"}
